# Apply DATA UPDATE: append new product rows (119-150) to the tracking sheet,
# and drop the now-obsolete trailing empty placeholder cells (K118:R118) that
# used to mark the end of the table (row 118 was previously the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old last row (118) had empty placeholder cells in K:R; since row 118 is
# no longer the last row of the table, those placeholders are removed.
$ws.Range("K118:R118").ClearContents()

# Helper: write a value as Text so numeric-looking codes/quantities (e.g. "40",
# "100", "1") are preserved as text cells instead of being auto-converted to
# numbers, matching the rest of the column. The style is reset back to
# "Normal" right after so we don't leave a stray text-format style behind.
function Set-TextCell($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

$rows = @(
    @{ Row=119; A='0TF26701'; B='BONNYHILL B VITAMIN-C SERUM'; C='VARIOS'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='40'; H='ML'; J='Revisado y Traducido' },
    @{ Row=120; A='0TF26702'; B='BONNYHILL FOOT CREAM 100ML'; C='VARIOS'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='100'; H='ML'; J='Revisado y Traducido' },
    @{ Row=121; A='0TS04122'; B='ZIAJA SET RUTINA AUTOBRONCEADOR'; C='TRAT.SOLAR'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='4'; H='UND'; J='Revisado y Traducido' },
    @{ Row=122; A='2BG01773'; B='ZIAJA CUPUAZU JABON CRISTALINO BAÑO & DUCHA 500ML'; C='TRAT.SOLAR'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=123; A='2BL02930'; B='ZIAJA CUPUAZU LOCION CORPORAL BRONCEADORA 300ML'; C='TRAT.SOLAR'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=124; A='0TN00932'; B='ZIAJA CUPUAZU EXFOLIANTE DE AZUCAR CRISTALINO 200M'; C='TRAT.SOLAR'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=125; A='0TS03888'; B='ZIAJA SOLAR CREMA FACIAL BRONCEADORA BRONZE 50ML'; C='TRAT.SOLAR'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=126; A='2CC02909'; B='ORIGINAL REMEDIES CHAMPU  5 PLANTAS 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=127; A='2CC05860'; B='ORIGINAL REMEDIES CHAMPU 5 PLANTAS 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=128; A='2CC05858'; B='ORIGINAL REMEDIES CHAMPU AGUA ARROZ 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=129; A='2CC02914'; B='ORIGINAL REMEDIES CHAMPU AGUA COCO & ALOE 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=130; A='2CC05853'; B='ORIGINAL REMEDIES CHAMPU AGUA COCO 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=131; A='2CC05606'; B='ORIGINAL REMEDIES CHAMPU AGUA DE ARROZ 300ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=132; A='2CC02904'; B='ORIGINAL REMEDIES CHAMPU AGUACATE & KARITE 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=133; A='2CC05854'; B='ORIGINAL REMEDIES CHAMPU AGUACATE 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=134; A='2CC02908'; B='ORIGINAL REMEDIES CHAMPU ARCILLA & LIMON 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=135; A='2CC05855'; B='ORIGINAL REMEDIES CHAMPU ARGAN 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=136; A='2CC05862'; B='ORIGINAL REMEDIES CHAMPU AVENA 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=137; A='2CC02911'; B='ORIGINAL REMEDIES CHAMPU AVENA DELICATESSE 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=138; A='2CC05857'; B='ORIGINAL REMEDIES CHAMPU CAMOMILA 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=139; A='2CC04848'; B='ORIGINAL REMEDIES CHAMPU CARBON 300ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=140; A='2CC05863'; B='ORIGINAL REMEDIES CHAMPU CARBON 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=141; A='2CC05019'; B='ORIGINAL REMEDIES CHAMPU CHARCOAL 250ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=142; A='2CC02903'; B='ORIGINAL REMEDIES CHAMPU ELIXIR ARGAN 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=143; A='2CC05861'; B='ORIGINAL REMEDIES CHAMPU LIMON 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=144; A='2CC02910'; B='ORIGINAL REMEDIES CHAMPU MIEL 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=145; A='2CC05859'; B='ORIGINAL REMEDIES CHAMPU MIEL 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=146; A='2CC05856'; B='ORIGINAL REMEDIES CHAMPU OLIVA 400ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=147; A='2CC02912'; B='ORIGINAL REMEDIES CHAMPU OLIVA MITICA 300 ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Revisado y Traducido' },
    @{ Row=148; A='2CC05918'; B='HERBAL CHAMPU COCO PACK 2X350ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='2'; H='UND'; J='Revisado y Traducido' },
    @{ Row=149; A='2CC05919'; B='HERBAL DESIRE ROSE PACK CH350ML+MASC500ML'; C='CABELLO CHAMPU'; D='No Tiene PT - TRADUZIDO'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; I='"8700216211086","8700216200547"'; J='Solo Revisión' },
    @{ Row=150; A='2CC05680'; B='HERBAL BIO RENEW CHAMPU ROSA 350ML'; C='CABELLO CHAMPU'; D='Tiene PT'; E='Tiene ES'; F='No Tiene IT - TRADOTTO'; G='1'; H='UND'; J='Solo Revisión' },
)

foreach ($row in $rows) {
    $r = $row.Row
    $addrA = "A" + $r
    $addrB = "B" + $r
    $addrC = "C" + $r
    $addrD = "D" + $r
    $addrE = "E" + $r
    $addrF = "F" + $r
    $addrG = "G" + $r
    $addrH = "H" + $r
    $addrI = "I" + $r
    $addrJ = "J" + $r
    Set-TextCell $ws $addrA $row.A
    Set-TextCell $ws $addrB $row.B
    Set-TextCell $ws $addrC $row.C
    Set-TextCell $ws $addrD $row.D
    Set-TextCell $ws $addrE $row.E
    Set-TextCell $ws $addrF $row.F
    Set-TextCell $ws $addrG $row.G
    Set-TextCell $ws $addrH $row.H
    if ($row.ContainsKey("I")) {
        Set-TextCell $ws $addrI $row.I
    }
    Set-TextCell $ws $addrJ $row.J
}

# Row 150 is now the last row of the table; mirror the same "end of table"
# placeholder pattern that row 118 used to have (an untouched/empty I cell
# plus empty K:R cells), so the used range/dimension extends through column R.
$ws.Range("I150").Style = "Normal"
$ws.Range("K150:R150").Style = "Normal"
